# Weekly fruit/vegetable price update: insert a new price record for
# "Bruselas (repollito)" at Vega Modelo de Temuco, shifting the existing
# rows 70-84 down to 71-85.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 70; this pushes rows 70..84 down to 71..85
# and extends the used range from A1:R84 to A1:R85.
$ws.Rows(70).Insert()

# Populate the newly inserted row 70 with this week's record.
$ws.Cells.Item(70, 1).Value = 10
$ws.Cells.Item(70, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(70, 3).Value = "La Araucanía"
$ws.Cells.Item(70, 4).Value = 44722
$ws.Cells.Item(70, 5).Value = 9
$ws.Cells.Item(70, 6).Value = 100112035
$ws.Cells.Item(70, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(70, 8).Value = "Sin especificar"
$ws.Cells.Item(70, 9).Value = "Primera"
$ws.Cells.Item(70, 10).Value = 20
$ws.Cells.Item(70, 11).Value = 28000
$ws.Cells.Item(70, 12).Value = 28000
$ws.Cells.Item(70, 13).Value = 28000
$ws.Cells.Item(70, 14).Value = "$/malla 10 kilos"
$ws.Cells.Item(70, 15).Value = "Región Metropolitana"
$ws.Cells.Item(70, 16).Value = 2800
$ws.Cells.Item(70, 17).Value = 10
$ws.Cells.Item(70, 18).Value = "Hortaliza"
